# Generate Report for Archive
#
# The e2e files "a8b14d5f-5967-415e-ab28-adf143a51eb5.md" (row 5) and
# "b17fc8e9-9142-4e40-8466-79fe1de02ba3.md" (row 6) moved back from
# "Ready for handoff" to "In Translation" on every sheet of the
# localization-status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("C6").Value = "In Translation"

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("C6").Value = "In Translation"
